$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (SST) - all values become 0
$ws.Range("B2").Value = 0
$ws.Range("B3").Value = 0
$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 0
$ws.Range("B7").Value = 0
$ws.Range("B8").Value = 0
$ws.Range("B9").Value = 0

# Column C (MSLP)
$ws.Range("C2").Value = -0.8127810846533238
$ws.Range("C3").Value = 0
$ws.Range("C4").Value = 0.7890902412313365
$ws.Range("C5").Value = 0.7432027668855248
$ws.Range("C6").Value = 0.7751198501047998
$ws.Range("C7").Value = 0
$ws.Range("C8").Value = -0.8201281856513709
$ws.Range("C9").Value = 0.8403027920443081

# Column D (Z500)
$ws.Range("D2").Value = 0.7512506727001143
$ws.Range("D3").Value = 0
$ws.Range("D4").Value = 0.7242029148177772
$ws.Range("D5").Value = 0.7176548363381191
$ws.Range("D6").Value = -0.7812359619930944
$ws.Range("D7").Value = 0
$ws.Range("D8").Value = 0.8481888557945662
$ws.Range("D9").Value = 0.827733620692871
